# Updates cryptos list cell values (Price and Volume(1h) columns)
# to match the latest scraped data, per commit message:
# "Updated cryptos list on Tue Nov 26 05:58:22 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'94.639.38"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").Value = "'3.426.38"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'238.59"
$ws.Range("E5").Value = "  -5.79%  "
$ws.Range("D6").Value = "'644.23"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.983"
$ws.Range("E10").Value = "  -5.89%  "
$ws.Range("D11").Value = "'3.425.45"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("D13").Value = "'42.19"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "'94.410.35"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").Value = "'4.064.92"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "'8.41"
$ws.Range("E18").Value = "  -6.99%  "
$ws.Range("D19").Value = "'3.425.69"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").Value = "'17.57"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").Value = "'11.72"
$ws.Range("E21").Value = "  +6.92%  "
$ws.Range("D22").Value = "'0.498"
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("D23").Value = "'501.03"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'3.24"
$ws.Range("E24").Value = "  -5.74%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("D27").Value = "'94.52"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "'12.02"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "'3.608.69"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").Value = "'11.80"
$ws.Range("E30").Value = "  +3.21%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'2.77"
$ws.Range("E32").Value = "  +7.02%  "
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'29.76"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'570.44"
$ws.Range("E38").Value = "  +6.75%  "
$ws.Range("D39").Value = "'7.71"
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").Value = "'0.905"
$ws.Range("E43").Value = "  +5.80%  "
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "'3.71"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "'5.69"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("D48").Value = "'3.35"
$ws.Range("D49").Value = "'0.0412"
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("E51").Value = "  -3.35%  "
